# Fix verse references in column A that had a spurious trailing "16"
# appended to them (e.g. "2 Corinthians 1:116" -> "2 Corinthians 1:1").
# Rows that were already correct (no trailing "16") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 212 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($null -ne $val -and $val -is [string] -and $val.EndsWith("16")) {
        $cell.Value = $val.Substring(0, $val.Length - 2)
    }
}
